# Penalty Reward System update (unfinished) - shifts forecast week dates
# forward by one week and updates MyForecast values, then refreshes the
# Summary sheet metrics derived from the Forecast Comparison data.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Forecast Comparison" --------------------------------------
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

# New Week_Start_Date values (col B) - shifted forward one week
$newDates = @(
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20",
    "2025-04-27"
)

# New MyForecast values (col D)
$newForecast = @(53, 56, 55, 50, 43, 41, 44, 33, 31, 32, 31, 46, 46, 31, 31, 30)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $row = $i + 2
    $ws1.Cells.Item($row, 2).NumberFormat = "@"
    $ws1.Cells.Item($row, 2).Value = $newDates[$i]
    $ws1.Cells.Item($row, 4).Value = $newForecast[$i]
}

# --- Sheet 2: "Summary" ---------------------------------------------------
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Cells.Item(2, 2).NumberFormat = "@"
$ws2.Cells.Item(2, 2).Value = "2023-02-12 to 2025-01-05"

$ws2.Cells.Item(4, 2).NumberFormat = "@"
$ws2.Cells.Item(4, 2).Value = "93"

$ws2.Cells.Item(8, 2).NumberFormat = "@"
$ws2.Cells.Item(8, 2).Value = "2510 units"

$ws2.Cells.Item(9, 2).NumberFormat = "@"
$ws2.Cells.Item(9, 2).Value = "654"

$ws2.Cells.Item(10, 2).NumberFormat = "@"
$ws2.Cells.Item(10, 2).Value = "375"

$ws2.Cells.Item(11, 2).NumberFormat = "@"
$ws2.Cells.Item(11, 2).Value = "214"

$ws2.Cells.Item(12, 2).NumberFormat = "@"
$ws2.Cells.Item(12, 2).Value = "56"

$ws2.Cells.Item(13, 2).NumberFormat = "@"
$ws2.Cells.Item(13, 2).Value = "2025-01-19"

$ws2.Cells.Item(14, 2).NumberFormat = "@"
$ws2.Cells.Item(14, 2).Value = "30"

$ws2.Cells.Item(15, 2).NumberFormat = "@"
$ws2.Cells.Item(15, 2).Value = "2025-04-27"
